$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 51-54: the order boundaries were manually solved, so D/E/F become
# plain (manually entered) numbers instead of formulas, and G becomes a
# simple "=E-D" range formula (matching rows 52-61's existing pattern).

$ws.Range("D51").Value = 3705.1
$ws.Range("E51").Value = 3767
$ws.Range("F51").Value = 3736.7
$ws.Range("G51").Formula = "=E51-D51"

$ws.Range("D52").Value = 3746.2
$ws.Range("E52").Value = 3808.5
$ws.Range("F52").Value = 3777.8
$ws.Range("G52").Formula = "=E52-D52"

$ws.Range("D53").Value = 3787.5
$ws.Range("E53").Value = 3850.8
$ws.Range("F53").Value = 3819.8
$ws.Range("G53").Formula = "=E53-D53"

$ws.Range("D54").Value = 3830.1
$ws.Range("E54").Value = 3894.1
$ws.Range("F54").Value = 3862.8
$ws.Range("G54").Formula = "=E54-D54"

# Mark these newly-solved orders as "Auto Done" in column K, matching
# the existing label already present on rows 55-61.
$ws.Range("K51").Value = "Auto Done"
$ws.Range("K52").Value = "Auto Done"
$ws.Range("K53").Value = "Auto Done"
$ws.Range("K54").Value = "Auto Done"

$wb.Application.Calculate()

$ws.Range("K52").Select()
